# Replace the leading word "Public" with "Code" in the Restoration
# Framework table row description, splitting it into its own run
# (matching how Word splits a run when only part of its text is
# retyped over a selection) while keeping identical run formatting.

$d = $word.ActiveDocument

$rng = $d.Content
$rng.Find.Execute("Public", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$rng.Text = "Code"

# Nudge a character-level property on just this range and revert it so
# the engine keeps "Code" as its own run instead of re-merging it back
# into the run that follows.
$rng.Bold = 1
$rng.Bold = 0
